$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "51.606.55"
$ws.Range("E2").Value = "  +0.78%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.986.20"
$ws.Range("E3").Value = "  +2.45%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.15%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "381.03"
$ws.Range("E5").Value = "  +4.51%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "106.37"
$ws.Range("E6").Value = "  +1.68%  "
$ws.Range("E7").Value = "  +0.96%  "
$ws.Range("E8").Value = "  -0.06%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.599"
$ws.Range("E9").Value = "  +1.52%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "37.56"
$ws.Range("E10").Value = "  +1.42%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0845"
$ws.Range("E12").Value = "  +1.25%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "18.72"
$ws.Range("E13").Value = "  +1.34%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "3.461.31"
$ws.Range("E14").Value = "  +2.39%  "
$ws.Range("E15").Value = "  +2.09%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.991.24"
$ws.Range("E16").Value = "  +2.49%  "
$ws.Range("E17").Value = "  +2.03%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "51.634.38"
$ws.Range("E18").Value = "  +0.74%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.39"
$ws.Range("E19").Value = "  +2.34%  "
$ws.Range("E20").Value = "  +2.35%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.05"
$ws.Range("E21").Value = "  -0.01%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.0₃0959"
$ws.Range("E22").Value = "  +1.44%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "69.41"
$ws.Range("E23").Value = "  +1.77%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "264.12"
$ws.Range("E24").Value = "  +1.59%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.80"
$ws.Range("E25").Value = "  +3.85%  "
$ws.Range("E26").Value = "  -2.10%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.28"
$ws.Range("E27").Value = "  +18.58%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.44"
$ws.Range("E28").Value = "  +2.36%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("B29").Value = "Dai"
$ws.Range("C29").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D29").Value = "1.00"
$ws.Range("E29").Value = "  -0.04%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("B30").Value = "EthereumClassic"
$ws.Range("C30").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D30").Value = "26.08"
$ws.Range("E30").Value = "  +0.52%  "
$ws.Range("E31").Value = "  +3.37%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "9.90"
$ws.Range("E32").Value = "  -0.53%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("B33").Value = "VeChain"
$ws.Range("C33").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D33").Value = "0.0467"
$ws.Range("E33").Value = "  +10.28%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("B34").Value = "InjectiveProtocol"
$ws.Range("C34").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D34").Value = "34.75"
$ws.Range("E34").Value = "  -1.46%  "
$ws.Range("E35").Value = "  -2.14%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "51.26"
$ws.Range("E36").Value = "  +1.21%  "
$ws.Range("E37").Value = "  +0.07%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.11"
$ws.Range("E38").Value = "  -1.21%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "17.50"
$ws.Range("E39").Value = "  +3.03%  "
$ws.Range("E40").Value = "  -6.90%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.85"
$ws.Range("E41").Value = "  -0.68%  "
$ws.Range("E42").Value = "  +2.49%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "123.70"
$ws.Range("E43").Value = "  +4.72%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "22.21"
$ws.Range("E44").Value = "  -1.45%  "
$ws.Range("E45").Value = "  -1.18%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.277"
$ws.Range("E46").Value = "  +17.59%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.050.18"
$ws.Range("E47").Value = "  -0.83%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.35"
$ws.Range("E48").Value = "  +3.71%  "
$ws.Range("E49").Value = "  +1.72%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0353"
$ws.Range("E50").Value = "  +10.72%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "5.20"
$ws.Range("E51").Value = "  +3.43%  "

Write-Host "Cryptos list updated"
